$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scalar fixes on existing rows ---
# Q58: 2 -> 0
$ws.Cells.Item(58, 17).Value = 0
# Q61: 1 -> 0
$ws.Cells.Item(61, 17).Value = 0
# O351: 0 -> 2
$ws.Cells.Item(351, 15).Value = 2
# R353 / R354 were empty inline-string placeholders -> become numeric 0
$ws.Cells.Item(353, 18).Value = 0
$ws.Cells.Item(354, 18).Value = 0

# --- New weekly rows 355-372 (columns A..Q; R stays blank like the prior tail rows) ---
$newRows = @(
  @(45474, 1787, 1877.400024414062, 1782.150024414062, 1865.050048828125, 1859.732788085938, 3891080, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
  @(45481, 1865.050048828125, 1886.199951171875, 1822, 1865.300048828125, 1859.982055664062, 2689613, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
  @(45488, 1867.150024414062, 1916, 1853.550048828125, 1881.400024414062, 1876.0361328125, 2967834, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
  @(45495, 1893, 1980.050048828125, 1856.849975585938, 1967.300048828125, 1961.691284179688, 5796257, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
  @(45502, 1969, 2031, 1922.849975585938, 1983.300048828125, 1977.645629882812, 4435562, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
  @(45509, 1959.5, 1996.449951171875, 1905.050048828125, 1965.949951171875, 1960.344970703125, 3380118, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
  @(45516, 1973.300048828125, 2046.800048828125, 1930, 2036.699951171875, 2030.893310546875, 3529540, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
  @(45523, 2043.449951171875, 2119.300048828125, 1993.050048828125, 2096.800048828125, 2090.822021484375, 4455568, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
  @(45530, 2108, 2183, 2080.199951171875, 2141.75, 2135.643798828125, 6099052, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
  @(45537, 2142, 2285, 2142, 2261.699951171875, 2255.251708984375, 4604456, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
  @(45544, 2261.550048828125, 2273.449951171875, 2105.050048828125, 2109.25, 2103.236572265625, 3731439, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
  @(45551, 2114, 2254.5, 2076.5, 2234.050048828125, 2227.6806640625, 5159895, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
  @(45558, 2240, 2301.89990234375, 2200.14990234375, 2238.14990234375, 2231.768798828125, 3583019, 2024, 9, 23, 0, 0, 0, 39, 1, 0, 0),
  @(45565, 2236.5, 2239.800048828125, 2116.89990234375, 2126.60009765625, 2120.537109375, 2818093, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
  @(45572, 2150, 2165.800048828125, 2060, 2067.35009765625, 2061.4560546875, 2587162, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
  @(45579, 2065, 2145, 1988.550048828125, 2020.949951171875, 2015.188232421875, 2804428, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
  @(45586, 2034, 2094.35009765625, 1916.800048828125, 1929.150024414062, 1923.650024414062, 4540429, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
  @(45593, 1911.599975585938, 1956.199951171875, 1870.599975585938, 1918.800048828125, 1918.800048828125, 4685852, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0)
)

$dateFormat = $ws.Cells.Item(354, 1).NumberFormat()

$r = 355
foreach ($row in $newRows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  # column A carries the date/time display format, like every other row
  $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
  $r++
}
